$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Values are written with a leading apostrophe to force text storage
# (matching the original inline-string cell type), then the cell style
# is reset to Normal so no numeric/text formatting is introduced.

$ws.Range("D2").Value = "'290.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-6.28%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'39.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.84%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.031"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.41%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07348"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.37%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'4.286"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.09%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.552"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-8.75%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9114"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.90%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1191"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-6.92%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1743"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-4.96%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.08682"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.07%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.04162"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.58%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.1052"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.03%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001275"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.48%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.005873"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.19%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.395"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.31%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D18").Value = "'0.3282"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.14%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'7.544"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.13%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1352"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.73%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.2885"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'6.09%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.03846"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.23%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.001273"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.56%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.003885"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-8.39%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'0.80%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0003728"
$ws.Range("D26").Style = "Normal"

$ws.Range("D38").Value = "'0.02325"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-8.96%"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.05011"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-5.65%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.007698"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.75%"
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'163.00%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1270"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.45%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.007374"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'11.08%"
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'-13.75%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.3137"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.37%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006518"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-3.83%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'15.46%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.004206"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'35.54%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.12%"
$ws.Range("E51").Style = "Normal"

